$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the neighboring "sum" header (G1) onto the new
# "Save" header (H1) so it matches the bold/centered/bordered look of the
# other column headers, then set the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for the "Save" column, row 2.
$ws.Range("H2").Value = 0
